# Applies the WeeklyPlans/9-17 edit:
#  - trims trailing space from the PRIORITY bullet
#  - adds Zipline / Quantopian / Pyalgotrade sub-bullets
#  - repurposes the old "Add config file..." paragraph into a new
#    "SMA - " bullet (inheriting the _GoBack bookmark)
#  - re-adds the "Add config file for hiding API key" bullet (spell-checked)
#  - extends "Best way to do this?" with "- .gitignore" (bookmark removed)
#  - splits "Grab data into csv file..." into spell-checked runs, adds
#    "Df to csv" sub-bullet
#  - splits the "Google finance..." bullet into spell/grammar-checked runs
#    and appends "VPN through colgate and then try" / "Talk to Aaron"
#
# Strategy: build the exact target OOXML for each touched/new paragraph
# and push it in via Range.InsertXML — this lets us reproduce the
# proofErr/bookmark/run-splitting exactly as the diff specifies.
# We walk the document bottom-to-top so earlier paragraph indices never
# shift under us.

$d = $word.ActiveDocument

function New-PkgXml($innerParagraphsXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerParagraphsXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---- Paragraph fragments (exact target OOXML) ----------------------------

$p_priority = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>PRIORITY – figure out how to compare strategies against baseline</w:t></w:r></w:p>'

$p_zipline_quantopian_pyalgotrade = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Zipline</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Quantopian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Pyalgotrade</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'

$p_sma = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">SMA - </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$p_addconfig = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file for hiding API key</w:t></w:r></w:p>'

$p_bestway = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Best way to do this?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>- .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gitignore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$p_grabdata = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Grab data into </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>csv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file – pointless to keep downloading </w:t></w:r></w:p>'

$p_dftocsv = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Df</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>csv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$p_googlefinance = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Google finance – should work but perhaps </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>bc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>colgate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wifi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?</w:t></w:r></w:p>'

$p_vpn_talktoaaron = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">VPN through </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>colgate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and then try</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Talk to Aaron</w:t></w:r></w:p>'

# ---- Apply bottom-to-top so untouched indices above stay valid -----------

# 12: "Google finance..." -> split runs; then append VPN/Talk-to-Aaron bullets
$para = $d.Paragraphs.Item(12)
$para.Range.InsertXML((New-PkgXml $p_googlefinance))
$para = $d.Paragraphs.Item(12)
$para.Range.InsertParagraphAfter()
$d.Paragraphs.Item(13).Range.InsertXML((New-PkgXml $p_vpn_talktoaaron))

# 6: "Grab data into csv file..." -> split runs; then add "Df to csv" bullet
$para = $d.Paragraphs.Item(6)
$para.Range.InsertXML((New-PkgXml $p_grabdata))
$para = $d.Paragraphs.Item(6)
$para.Range.InsertParagraphAfter()
$d.Paragraphs.Item(7).Range.InsertXML((New-PkgXml $p_dftocsv))

# 5: "Best way to do this?" -> extend text, drop bookmark
$para = $d.Paragraphs.Item(5)
$para.Range.InsertXML((New-PkgXml $p_bestway))

# 4: "Add config file for hiding API key" -> becomes "SMA - " (+bookmark);
#    then re-insert the (spell-checked) "Add config file..." bullet after it
$para = $d.Paragraphs.Item(4)
$para.Range.InsertXML((New-PkgXml $p_sma))
$para = $d.Paragraphs.Item(4)
$para.Range.InsertParagraphAfter()
$d.Paragraphs.Item(5).Range.InsertXML((New-PkgXml $p_addconfig))

# 3: "PRIORITY..." -> trim trailing space; then add Zipline/Quantopian/Pyalgotrade
$para = $d.Paragraphs.Item(3)
$para.Range.InsertXML((New-PkgXml $p_priority))
$para = $d.Paragraphs.Item(3)
$para.Range.InsertParagraphAfter()
$d.Paragraphs.Item(4).Range.InsertXML((New-PkgXml $p_zipline_quantopian_pyalgotrade))

Write-Host "Final paragraph count:" $d.Paragraphs.Count
